$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("B4").Value = "SardanaRip"
$ws.Range("C4").Value = $null
$ws.Range("D4").Value = $null
$ws.Range("E4").Value = $null
$ws.Range("F4").Value = $null
$ws.Range("G4").Value = $null
$ws.Range("H4").Value = "Olaaa"
$ws.Range("I4").Value = "Olaaa"

# Row 5 - clear all cells
$ws.Range("A5:I5").Value = $null
